# Update the "想去人数" (want-to-go count) figures in column F for the
# sheets "展览" (sheet1) and "全部类型" (sheet4), reflecting a refreshed
# data scrape (gh-pages output regenerated at commit 456a3b4).
#
# Sheets "演出" and "本地生活" are unaffected by this refresh.

$wb = $excel.ActiveWorkbook

# row -> new F value, for worksheet "展览"
$exhibitionUpdates = @{
    2  = 135
    3  = 340
    4  = 429
    5  = 1738
    6  = 87
    7  = 2182
    9  = 285
    11 = 4931
    12 = 11
    17 = 185
    18 = 36
    19 = 22
    21 = 3884
    22 = 712
    23 = 669
    26 = 105
    27 = 120
    31 = 579
    33 = 28
    34 = 951
    35 = 2462
    36 = 426
}

# row -> new F value, for worksheet "全部类型"
$allTypesUpdates = @{
    2  = 135
    3  = 340
    4  = 429
    5  = 1738
    6  = 87
    7  = 2182
    9  = 285
    11 = 4931
    12 = 11
    17 = 185
    18 = 36
    19 = 22
    21 = 3884
    22 = 712
    23 = 669
    26 = 105
    27 = 120
    31 = 579
    34 = 28
    35 = 951
    36 = 2462
    37 = 426
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Cells.Item($row, 6).Value = $exhibitionUpdates[$row]
}

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Cells.Item($row, 6).Value = $allTypesUpdates[$row]
}
